$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "chronických onemocnění (0, 1- 2; 3-5; 6+)",
    $true,   # MatchCase
    $false,  # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap
    $false,  # Format
    "chronických onemocnění (0, 1- 2; 3-5; 6+)",
    2        # Replace: wdReplaceAll
)
